$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 302.78262
$ws.Range("I33").Value = 300.2857
$ws.Range("J33").Value = 306.66666
$ws.Range("K33").Value = 300.2857
$ws.Range("L33").Value = 306.66666
$ws.Range("M33").Value = -71.28570000000002
$ws.Range("N33").Value = -764.66666
$ws.Range("H55").Value = 70.90000000000001
$ws.Range("I55").Value = 70.90000000000001
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 70.90000000000001
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 143.1
$ws.Range("H62").Value = 3709.7
$ws.Range("I62").Value = 3059.4
$ws.Range("J62").Value = 4360
$ws.Range("K62").Value = 3059.4
$ws.Range("L62").Value = 4360
$ws.Range("M62").Value = -2435.4
$ws.Range("N62").Value = -5608
$ws.Range("H65").Value = 3709.7
$ws.Range("I65").Value = 3059.4
$ws.Range("J65").Value = 4360
$ws.Range("K65").Value = 15297
$ws.Range("L65").Value = 21800
$ws.Range("M65").Value = -12177
$ws.Range("N65").Value = -28040
$ws.Range("H98").Value = 1645
$ws.Range("I98").Value = 1338.125
$ws.Range("J98").Value = 4100
$ws.Range("K98").Value = 1338.125
$ws.Range("L98").Value = 4100
$ws.Range("M98").Value = 159.875
$ws.Range("N98").Value = -7096
$ws.Range("H122").Value = 1645
$ws.Range("I122").Value = 1338.125
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 4014.375
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -1564.375
$ws.Range("N122").Value = -17200
$ws.Range("H134").Value = 62222.855
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 62222.855
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 62222.855
$ws.Range("N134").Value = -72362.85500000001
$ws.Range("H135").Value = 166670450
$ws.Range("I135").Value = 71432720
$ws.Range("J135").Value = 500002500
$ws.Range("K135").Value = 642894480
$ws.Range("L135").Value = 4500022500
$ws.Range("M135").Value = -642891945
$ws.Range("N135").Value = -4500027570
$ws.Range("H136").Value = 74782.71000000001
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 74782.71000000001
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 74782.71000000001
$ws.Range("N136").Value = -84982.71000000001
$ws.Range("H137").Value = 3144.2292
$ws.Range("I137").Value = 1704.16
$ws.Range("J137").Value = 4709.522
$ws.Range("K137").Value = 5112.48
$ws.Range("L137").Value = 14128.566
$ws.Range("M137").Value = -2562.48
$ws.Range("N137").Value = -19228.566
$ws.Range("H138").Value = 3671.3098
$ws.Range("I138").Value = 1453.5333
$ws.Range("J138").Value = 4265.357
$ws.Range("K138").Value = 4360.5999
$ws.Range("L138").Value = 12796.071
$ws.Range("M138").Value = 779.4000999999998
$ws.Range("N138").Value = -23076.071

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 20000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 20000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -20644
$ws.Range("H32").Value = 5740.86
$ws.Range("I32").Value = 3769.7012
$ws.Range("J32").Value = 18932.46
$ws.Range("K32").Value = 3769.7012
$ws.Range("L32").Value = 18932.46
$ws.Range("M32").Value = -3482.7012
$ws.Range("N32").Value = -19506.46
$ws.Range("H74").Value = 5549.1665
$ws.Range("I74").Value = 3014
$ws.Range("J74").Value = 14422.25
$ws.Range("K74").Value = 3014
$ws.Range("L74").Value = 14422.25
$ws.Range("M74").Value = -2140
$ws.Range("N74").Value = -16170.25
$ws.Range("H77").Value = 5549.1665
$ws.Range("I77").Value = 3014
$ws.Range("J77").Value = 14422.25
$ws.Range("K77").Value = 15070
$ws.Range("L77").Value = 72111.25
$ws.Range("M77").Value = -10702
$ws.Range("N77").Value = -80847.25
$ws.Range("H122").Value = 7355506.5
$ws.Range("I122").Value = 4437.3335
$ws.Range("J122").Value = 8930736
$ws.Range("K122").Value = 13312.0005
$ws.Range("L122").Value = 26792208
$ws.Range("M122").Value = -10862.0005
$ws.Range("N122").Value = -26797108
$ws.Range("H132").Value = 6568.6665
$ws.Range("I132").Value = 1713.9
$ws.Range("J132").Value = 8679.434999999999
$ws.Range("K132").Value = 5141.700000000001
$ws.Range("L132").Value = 26038.305
$ws.Range("M132").Value = -2611.700000000001
$ws.Range("N132").Value = -31098.305

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 750
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = -360
$ws.Range("N11").Value = -1280
$ws.Range("H37").Value = 10100.5
$ws.Range("I37").Value = 201
$ws.Range("J37").Value = 20000
$ws.Range("K37").Value = 201
$ws.Range("L37").Value = 20000
$ws.Range("M37").Value = -64
$ws.Range("N37").Value = -20274
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 2102.2559
$ws.Range("I86").Value = 2157.6316
$ws.Range("J86").Value = 1681.4
$ws.Range("K86").Value = 2157.6316
$ws.Range("L86").Value = 1681.4
$ws.Range("M86").Value = -1034.6316
$ws.Range("N86").Value = -3927.4
$ws.Range("H89").Value = 2102.2559
$ws.Range("I89").Value = 2157.6316
$ws.Range("J89").Value = 1681.4
$ws.Range("K89").Value = 10788.158
$ws.Range("L89").Value = 8407
$ws.Range("M89").Value = -5172.158000000001
$ws.Range("N89").Value = -19639
$ws.Range("H134").Value = 4967.3335
$ws.Range("I134").Value = 4641.364
$ws.Range("J134").Value = 5863.75
$ws.Range("K134").Value = 13924.092
$ws.Range("L134").Value = 17591.25
$ws.Range("M134").Value = -11389.092
$ws.Range("N134").Value = -22661.25
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 43225.715
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 43225.715
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 43225.715
$ws.Range("N141").Value = -53585.715

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2555.7627
$ws.Range("I31").Value = 1868.6154
$ws.Range("J31").Value = 3895.7
$ws.Range("K31").Value = 1868.6154
$ws.Range("L31").Value = 3895.7
$ws.Range("M31").Value = -1573.6154
$ws.Range("N31").Value = -4485.7
$ws.Range("H34").Value = 2555.7627
$ws.Range("I34").Value = 1868.6154
$ws.Range("J34").Value = 3895.7
$ws.Range("K34").Value = 1868.6154
$ws.Range("L34").Value = 3895.7
$ws.Range("M34").Value = -1666.6154
$ws.Range("N34").Value = -4299.7
$ws.Range("H58").Value = 2276041
$ws.Range("I58").Value = 3790539
$ws.Range("J58").Value = 4293.6875
$ws.Range("K58").Value = 3790539
$ws.Range("L58").Value = 4293.6875
$ws.Range("M58").Value = -3790336
$ws.Range("N58").Value = -4699.6875
$ws.Range("H132").Value = 3748.3333
$ws.Range("I132").Value = 2745.875
$ws.Range("J132").Value = 5753.25
$ws.Range("K132").Value = 8237.625
$ws.Range("L132").Value = 17259.75
$ws.Range("M132").Value = -5707.625
$ws.Range("N132").Value = -22319.75
$ws.Range("H134").Value = 3589.6
$ws.Range("I134").Value = 2187.0557
$ws.Range("J134").Value = 4524.6294
$ws.Range("K134").Value = 6561.1671
$ws.Range("L134").Value = 13573.8882
$ws.Range("M134").Value = -4026.1671
$ws.Range("N134").Value = -18643.8882
$ws.Range("H136").Value = 2276041
$ws.Range("I136").Value = 3790539
$ws.Range("J136").Value = 4293.6875
$ws.Range("K136").Value = 11371617
$ws.Range("L136").Value = 12881.0625
$ws.Range("M136").Value = -11369067
$ws.Range("N136").Value = -17981.0625

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 27778028
$ws.Range("I12").Value = 83333570
$ws.Range("J12").Value = 257.41666
$ws.Range("K12").Value = 250000710
$ws.Range("L12").Value = 772.2499799999999
$ws.Range("M12").Value = -250000537
$ws.Range("N12").Value = -1118.24998
$ws.Range("H38").Value = 71.22221999999999
$ws.Range("I38").Value = 29.375
$ws.Range("J38").Value = 104.7
$ws.Range("K38").Value = 88.125
$ws.Range("L38").Value = 314.1
$ws.Range("M38").Value = 258.875
$ws.Range("N38").Value = -1008.1
$ws.Range("H63").Value = 2998.2354
$ws.Range("I63").Value = 2489.1428
$ws.Range("J63").Value = 3354.6
$ws.Range("K63").Value = 7467.428400000001
$ws.Range("L63").Value = 10063.8
$ws.Range("M63").Value = -6718.428400000001
$ws.Range("N63").Value = -11561.8
$ws.Range("H66").Value = 2998.2354
$ws.Range("I66").Value = 2489.1428
$ws.Range("J66").Value = 3354.6
$ws.Range("K66").Value = 22402.2852
$ws.Range("L66").Value = 30191.4
$ws.Range("M66").Value = -18658.2852
$ws.Range("N66").Value = -37679.39999999999
$ws.Range("H110").Value = 2709.7026
$ws.Range("I110").Value = 1413.5
$ws.Range("J110").Value = 2783.7715
$ws.Range("K110").Value = 4240.5
$ws.Range("L110").Value = 8351.3145
$ws.Range("M110").Value = -150.5
$ws.Range("N110").Value = -16531.3145
$ws.Range("H122").Value = 742
$ws.Range("I122").Value = 385.29413
$ws.Range("J122").Value = 1121
$ws.Range("K122").Value = 3467.64717
$ws.Range("L122").Value = 10089
$ws.Range("M122").Value = -1017.64717
$ws.Range("N122").Value = -14989

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 20998
$ws.Range("I17").Value = 8000
$ws.Range("J17").Value = 40495
$ws.Range("K17").Value = 8000
$ws.Range("L17").Value = 40495
$ws.Range("M17").Value = -7832
$ws.Range("N17").Value = -40831
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").ClearContents()
$ws.Range("H132").Value = 9526
$ws.Range("I132").Value = 12177
$ws.Range("J132").Value = 4224
$ws.Range("K132").Value = 36531
$ws.Range("L132").Value = 12672
$ws.Range("M132").Value = -34001
$ws.Range("N132").Value = -17732
$ws.Range("H137").Value = 49800
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 49800
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 49800
$ws.Range("N137").Value = -60000

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 11400
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 22300
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 22300
$ws.Range("M4").Value = -387
$ws.Range("N4").Value = -22526
$ws.Range("H20").Value = 11400
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 11400
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 11400
$ws.Range("N20").Value = -11852
$ws.Range("H25").Value = 32252
$ws.Range("I25").Value = 8500
$ws.Range("J25").Value = 56004
$ws.Range("K25").Value = 8500
$ws.Range("L25").Value = 56004
$ws.Range("M25").Value = -8270
$ws.Range("N25").Value = -56464
$ws.Range("H28").Value = 11400
$ws.Range("I28").Value = 500
$ws.Range("J28").Value = 22300
$ws.Range("K28").Value = 500
$ws.Range("L28").Value = 22300
$ws.Range("M28").Value = -268
$ws.Range("N28").Value = -22764
$ws.Range("H37").Value = 11400
$ws.Range("I37").Value = 500
$ws.Range("J37").Value = 22300
$ws.Range("K37").Value = 500
$ws.Range("L37").Value = 22300
$ws.Range("M37").Value = -393
$ws.Range("N37").Value = -22514
$ws.Range("H132").Value = 2765.1738
$ws.Range("I132").Value = 2128.0833
$ws.Range("J132").Value = 3460.182
$ws.Range("K132").Value = 6384.249899999999
$ws.Range("L132").Value = 10380.546
$ws.Range("M132").Value = -3854.249899999999
$ws.Range("N132").Value = -15440.546

